$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the date stamp "THU Dec 14" / " 10:16:53 PST 2017" was split
# across two runs; collapse it back into a single run with identical
# text by doing a no-op Find/Replace over the whole phrase.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "THU Dec 14 10:16:53 PST 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "THU Dec 14 10:16:53 PST 2017", 2)

# ---------------------------------------------------------------------
# Change 2: append a brand-new "chick-in" record (MON Dec 18, TSV /
# CARROT) right after the last "Amount balance" line, before the
# trailing blank paragraphs at the end of the document.
# ---------------------------------------------------------------------

# Locate the last paragraph that starts with "Amount balance" -- that is
# the bold "Amount balance ... - 834.0" line the new block must follow.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "^Amount balance") {
        $anchorIndex = $i
    }
}

# The paragraph right after the anchor is the first of the pre-existing
# trailing blank paragraphs; all new paragraphs are inserted just before
# it (i.e. right after the anchor).
$tailIndex = $anchorIndex + 1
$tailPara = $d.Paragraphs.Item($tailIndex)
$tailStart = $tailPara.Range.Start

# Carve out 9 new (still empty) paragraphs immediately before the tail
# paragraph, in document order.
for ($n = 1; $n -le 9; $n++) {
    $insPt = $d.Range($tailStart, $tailStart)
    $insPt.InsertParagraphBefore()
}

# Index of the first newly-created paragraph, and the rest computed
# eagerly into their own variables (avoids passing "$var (expr)" as
# adjacent command arguments, which this interpreter mis-parses).
$idx1 = $anchorIndex + 1
$idx2 = $idx1 + 1
$idx3 = $idx1 + 2
$idx4 = $idx1 + 3
$idx5 = $idx1 + 4
$idx6 = $idx1 + 5
$idx7 = $idx1 + 6
$idx8 = $idx1 + 7
$idx9 = $idx1 + 8

# Helper: append one run of text to the end of a paragraph (just before
# its paragraph mark) with the given formatting.
function Fill-Para($doc, $index, $txt, $bold, $color) {
    $pEnd = $doc.Paragraphs.Item($index).Range.End
    $insPos = $pEnd - 1
    $r = $doc.Range($insPos, $insPos)
    $r.Font.Name = "Courier New"
    if ($bold) { $r.Font.Bold = $true }
    if ($color -ne $null) { $r.Font.Color = $color }
    $r.InsertAfter($txt)
}

# Helper: strip the stray empty run Word's paragraph-insert leaves
# behind on an otherwise-empty paragraph, while preserving whatever
# paragraph-mark formatting (bold) was requested.
function Clean-EmptyPara($doc, $index, $bold) {
    $r = $doc.Paragraphs.Item($index).Range
    $start = $r.Start
    $ins = $doc.Range($start, $start)
    $ins.InsertAfter("X")
    $cEnd = $start + 1
    $c = $doc.Range($start, $cEnd)
    $c.Font.Name = "Courier New"
    if ($bold) { $c.Font.Bold = $true }
    $c.Delete()
}

$tab = [char]9

# 1) blank bold paragraph
Clean-EmptyPara $d $idx1 $true

# 2) "MON Dec 18" / " 10:30:04 PST 2017"
Fill-Para $d $idx2 "MON Dec 18" $false $null
Fill-Para $d $idx2 " 10:30:04 PST 2017" $false $null

# 3) "Person Name" <tab><tab><tab><tab>"- TSV"
Fill-Para $d $idx3 "Person Name" $false $null
Fill-Para $d $idx3 $tab $false $null
Fill-Para $d $idx3 $tab $false $null
Fill-Para $d $idx3 $tab $false $null
$tsvTail = $tab + "- TSV"
Fill-Para $d $idx3 $tsvTail $false $null

# 4) dashed separator line
Fill-Para $d $idx4 "---------------------------------------------------------------" $false $null

# 5) "Item Name" <tab><tab><tab><tab>"- CARROT"
Fill-Para $d $idx5 "Item Name" $false $null
Fill-Para $d $idx5 $tab $false $null
Fill-Para $d $idx5 $tab $false $null
Fill-Para $d $idx5 $tab $false $null
$carrotTail = $tab + "- CARROT"
Fill-Para $d $idx5 $carrotTail $false $null

# 6) "Amount Received" <tab><tab><tab>"- 834" (all red)
Fill-Para $d $idx6 "Amount Received" $false 255
Fill-Para $d $idx6 $tab $false 255
Fill-Para $d $idx6 $tab $false 255
$amt834Tail = $tab + "- 834"
Fill-Para $d $idx6 $amt834Tail $false 255

# 7) "Amount Received mode" <tab><tab>"- CASH AND CLEARD"
Fill-Para $d $idx7 "Amount Received mode" $false $null
Fill-Para $d $idx7 $tab $false $null
$cashTail = $tab + "- CASH AND CLEARD"
Fill-Para $d $idx7 $cashTail $false $null

# 8) blank paragraph (plain)
Clean-EmptyPara $d $idx8 $false

# 9) blank bold paragraph
Clean-EmptyPara $d $idx9 $true
